# Auto-applied data refresh: update H:N pricing/profit columns per scheduled runner pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value2 = 397.375
$ws.Range("I4").Value2 = 196.5
$ws.Range("K4").Value2 = 196.5
$ws.Range("M4").Value2 = -82.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value2 = 272.125
$ws.Range("I6").Value2 = 302.42856
$ws.Range("J6").Value2 = 60
$ws.Range("K6").Value2 = 907.28568
$ws.Range("L6").Value2 = 180
$ws.Range("M6").Value2 = -795.28568
$ws.Range("N6").Value2 = -404

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value2 = 3855.7144
$ws.Range("I18").Value2 = 3998.3333
$ws.Range("K18").Value2 = 3998.3333
$ws.Range("M18").Value2 = -3714.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value2 = 2950.75
$ws.Range("I28").Value2 = 1318.8
$ws.Range("K28").Value2 = 1318.8
$ws.Range("M28").Value2 = -833.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value2 = 1731.75
$ws.Range("I39").Value2 = 823.7692
$ws.Range("K39").Value2 = 2471.3076
$ws.Range("M39").Value2 = -2175.3076

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value2 = 460.55554
$ws.Range("I41").Value2 = 612.4
$ws.Range("J41").Value2 = 270.75
$ws.Range("K41").Value2 = 612.4
$ws.Range("L41").Value2 = 270.75
$ws.Range("M41").Value2 = -172.4
$ws.Range("N41").Value2 = -1150.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value2 = 3950.5
$ws.Range("I42").Value2 = 349
$ws.Range("K42").Value2 = 1047
$ws.Range("M42").Value2 = -817

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value2 = 3449.5334
$ws.Range("I62").Value2 = 3140.7273
$ws.Range("K62").Value2 = 3140.7273
$ws.Range("M62").Value2 = -2516.7273

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value2 = 3449.5334
$ws.Range("I65").Value2 = 3140.7273
$ws.Range("K65").Value2 = 15703.6365
$ws.Range("M65").Value2 = -12583.6365

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value2 = 7160.3076
$ws.Range("I113").Value2 = 7432.0713
$ws.Range("K113").Value2 = 7432.0713
$ws.Range("M113").Value2 = -4178.0713

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value2 = 1415.5264
$ws.Range("J137").Value2 = 1437.7916
$ws.Range("L137").Value2 = 4313.3748
$ws.Range("N137").Value2 = -9413.3748

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value2 = 3067.2114
$ws.Range("I138").Value2 = 2210.9167
$ws.Range("J138").Value2 = 3324.1
$ws.Range("K138").Value2 = 6632.750100000001
$ws.Range("L138").Value2 = 9972.299999999999
$ws.Range("M138").Value2 = -1492.750100000001
$ws.Range("N138").Value2 = -20252.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 2370.647
$ws.Range("I2").Value2 = 1653.4
$ws.Range("J2").Value2 = 7750
$ws.Range("K2").Value2 = 1653.4
$ws.Range("L2").Value2 = 7750
$ws.Range("M2").Value2 = -1540.4
$ws.Range("N2").Value2 = -7976

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value2 = 2440.9092
$ws.Range("I45").Value2 = 1500.2307
$ws.Range("J45").Value2 = 3799.6667
$ws.Range("K45").Value2 = 1500.2307
$ws.Range("L45").Value2 = 3799.6667
$ws.Range("M45").Value2 = -1123.2307
$ws.Range("N45").Value2 = -4553.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 3082.3333
$ws.Range("I61").Value2 = 2248.5
$ws.Range("K61").Value2 = 2248.5
$ws.Range("M61").Value2 = -2036.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value2 = 0
$ws.Range("J62").Value2 = 0
$ws.Range("L62").Value2 = 0
$ws.Range("N62").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H65").Value2 = 0
$ws.Range("J65").Value2 = 0
$ws.Range("L65").Value2 = 0
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value2 = 2370.647
$ws.Range("I116").Value2 = 1653.4
$ws.Range("J116").Value2 = 7750
$ws.Range("K116").Value2 = 1653.4
$ws.Range("L116").Value2 = 7750
$ws.Range("M116").Value2 = 640.5999999999999
$ws.Range("N116").Value2 = -12338

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value2 = 3082.3333
$ws.Range("I136").Value2 = 2248.5
$ws.Range("K136").Value2 = 6745.5
$ws.Range("M136").Value2 = -4195.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 2370.647
$ws.Range("I3").Value2 = 1653.4
$ws.Range("J3").Value2 = 7750
$ws.Range("K3").Value2 = 1653.4
$ws.Range("L3").Value2 = 7750
$ws.Range("M3").Value2 = -1539.4
$ws.Range("N3").Value2 = -7978

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value2 = 178.5
$ws.Range("I7").Value2 = 190.2
$ws.Range("K7").Value2 = 190.2
$ws.Range("M7").Value2 = -77.19999999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value2 = 2444.1667
$ws.Range("I99").Value2 = 1932.3334
$ws.Range("K99").Value2 = 1932.3334
$ws.Range("M99").Value2 = -434.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value2 = 2444.1667
$ws.Range("I126").Value2 = 1932.3334
$ws.Range("K126").Value2 = 5797.0002
$ws.Range("M126").Value2 = -3327.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value2 = 2518
$ws.Range("I132").Value2 = 2457.7917
$ws.Range("J132").Value2 = 2999.6667
$ws.Range("K132").Value2 = 7373.375100000001
$ws.Range("L132").Value2 = 8999.000100000001
$ws.Range("M132").Value2 = -4843.375100000001
$ws.Range("N132").Value2 = -14059.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I4").Value2 = 586688.3
$ws.Range("K4").Value2 = 1760064.9
$ws.Range("M4").Value2 = -1759952.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 879.7143
$ws.Range("J5").Value2 = 978.7
$ws.Range("L5").Value2 = 2936.1
$ws.Range("N5").Value2 = -3160.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value2 = 874.6667
$ws.Range("I69").Value2 = 886.5
$ws.Range("J69").Value2 = 780
$ws.Range("K69").Value2 = 2659.5
$ws.Range("L69").Value2 = 2340
$ws.Range("M69").Value2 = -1848.5
$ws.Range("N69").Value2 = -3962

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value2 = 874.6667
$ws.Range("I72").Value2 = 886.5
$ws.Range("J72").Value2 = 780
$ws.Range("K72").Value2 = 7978.5
$ws.Range("L72").Value2 = 7020
$ws.Range("M72").Value2 = -3922.5
$ws.Range("N72").Value2 = -15132

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value2 = 2776.15
$ws.Range("J107").Value2 = 2973.1667
$ws.Range("L107").Value2 = 8919.500100000001
$ws.Range("N107").Value2 = -12759.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value2 = 5773.364
$ws.Range("I117").Value2 = 1860.6666
$ws.Range("J117").Value2 = 7240.625
$ws.Range("K117").Value2 = 5581.9998
$ws.Range("L117").Value2 = 21721.875
$ws.Range("M117").Value2 = -2139.9998
$ws.Range("N117").Value2 = -28605.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value2 = 2217.111
$ws.Range("I118").Value2 = 2217.111
$ws.Range("K118").Value2 = 6651.333
$ws.Range("M118").Value2 = -5408.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value2 = 879.7143
$ws.Range("J135").Value2 = 978.7
$ws.Range("L135").Value2 = 8808.300000000001
$ws.Range("N135").Value2 = -13878.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value2 = 5499.8335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 4582.2383
$ws.Range("I70").Value2 = 4362.2
$ws.Range("K70").Value2 = 4362.2
$ws.Range("M70").Value2 = -4092.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value2 = 4582.2383
$ws.Range("I73").Value2 = 4362.2
$ws.Range("K73").Value2 = 4362.2
$ws.Range("M73").Value2 = -3426.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 15981
$ws.Range("I80").Value2 = 12999.6
$ws.Range("J80").Value2 = 20950
$ws.Range("K80").Value2 = 12999.6
$ws.Range("L80").Value2 = 20950
$ws.Range("M80").Value2 = -12001.6
$ws.Range("N80").Value2 = -22946

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value2 = 15981
$ws.Range("I83").Value2 = 12999.6
$ws.Range("J83").Value2 = 20950
$ws.Range("K83").Value2 = 64998
$ws.Range("L83").Value2 = 104750
$ws.Range("M83").Value2 = -60006
$ws.Range("N83").Value2 = -114734

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value2 = 70734.97
$ws.Range("J113").Value2 = 87478.75
$ws.Range("L113").Value2 = 87478.75
$ws.Range("N113").Value2 = -91818.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value2 = 46320.61
$ws.Range("I132").Value2 = 55574.05
$ws.Range("K132").Value2 = 166722.15
$ws.Range("M132").Value2 = -164192.15

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 5532.5625
$ws.Range("I7").Value2 = 4835.1113
$ws.Range("K7").Value2 = 4835.1113
$ws.Range("M7").Value2 = -4723.1113

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 8000
$ws.Range("I16").Value2 = 8000
$ws.Range("K16").Value2 = 8000
$ws.Range("M16").Value2 = -7830

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value2 = 735199
$ws.Range("I88").Value2 = 11996.25
$ws.Range("J88").Value2 = 998181.8
$ws.Range("K88").Value2 = 11996.25
$ws.Range("L88").Value2 = 998181.8
$ws.Range("M88").Value2 = -11568.25
$ws.Range("N88").Value2 = -999037.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H91").Value2 = 735199
$ws.Range("I91").Value2 = 11996.25
$ws.Range("J91").Value2 = 998181.8
$ws.Range("K91").Value2 = 11996.25
$ws.Range("L91").Value2 = 998181.8
$ws.Range("M91").Value2 = -10514.25
$ws.Range("N91").Value2 = -1001145.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value2 = 5532.5625
$ws.Range("I126").Value2 = 4835.1113
$ws.Range("K126").Value2 = 14505.3339
$ws.Range("M126").Value2 = -12035.3339
